# grade_book.xlsx edit
#
# The gradebook's "Final Grade" column (F) was left blank for every
# student. The author filled it in with each student's Final Paper score
# (column E) -- copying both the values and the existing cell formatting
# from column E so the new column F cells look just like the other score
# columns (centered, same font) instead of the empty placeholder style.
#
# Finally, the active selection in the sheet is left on G2 (the next
# empty column), matching where the user clicked after finishing the
# paste.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRows = $ws.Range("E2:E25")
$target   = $ws.Range("F2:F25")

# Bring over the values first ...
$dataRows.Copy()
$target.PasteSpecial(-4163)   # xlPasteValues

# ... then the formatting, so F2:F25 matches the styling already used by
# B2:E25 (centered, non-bold data font) instead of the blank-cell style.
$dataRows.Copy()
$target.PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# Leave the selection on G2, same as the finished workbook.
$ws.Range("G2").Select()
